# Issue List / Issues.xlsx update
# - Issue_002 (row 3) is now "Solved", with a new comment in column F.
# - A new Issue_003 (row 4) is logged as "In work".
# - Column F is widened to comfortably fit the long comment text.
# - Selection cursor ends up on F9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new Issue_003 row first (A4:C4) so the shared-string table picks up
# "Issue_003" / "Right click plugin isn't appearing" before "Solved" / the
# long comment text, matching the original authoring order.
$ws.Range("A4").Value = "Issue_003"
$ws.Range("B4").Value = "Right click plugin isn't appearing"
$ws.Range("C4").Value = 4

# Update Issue_002's status and add the resolution comment.
$ws.Range("D3").Value = "Solved"
$ws.Range("F3").Value = "added in the code to translate the wire type to the new one required by the ""Inline or Forked.vi"". Now appears to be working."

# Finish the new Issue_003 row - status reuses the existing "In work" string.
$ws.Range("D4").Value = "In work"

# Widen column F so the long comment text is readable.
$ws.Columns.Item(6).ColumnWidth = 146.33

# Leave the selection where the author left it after typing the comment.
$ws.Range("F9").Select()
